# feat: add parameters to check box tests
#
# Adds a new "CheckBox" worksheet (after the existing "userPositive" sheet)
# that holds the test parameters used by the CheckBox tests.

$wb = $excel.ActiveWorkbook

# Use the first sheet as a template so the new sheet naturally inherits the
# same page setup (paper size / orientation) that the other data sheets use.
$template = $wb.Worksheets.Item(1)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $lastSheet)

$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "CheckBox"
$ws.Cells.Clear()

# Header row
$ws.Range("A1").Value = "Box"
$ws.Range("B1").Value = "Results"
$ws.Range("A1:B1").Font.Bold = $true

# Data rows (written in this specific order so the shared-string table ends
# up populated in the same sequence the workbook was originally authored in)
$ws.Range("A4").Value = "Desktop React"
$ws.Range("B4").Value = "desktop notes commands react"

$ws.Range("A5").Value = "Classified"
$ws.Range("B5").Value = "classified"

$ws.Range("B2").Value = "home desktop notes commands documents workspace react angular veu office public private classified general downloads wordFile excelFile"
$ws.Range("A2").Value = "Home"

$ws.Range("A3").Value = "Home WorkSpace Public"
$ws.Range("B3").Value = "desktop notes commands private classified general downloads wordFile excelFile"

# Column widths roughly matching the authored workbook (A ~ 22.9, B = 29)
$ws.Columns.Item(1).ColumnWidth = 22.0
$ws.Columns.Item(2).ColumnWidth = 28.17
